# "style a new body" -- the original single journal paragraph (9.24/9.25's
# entry about js) is kept, but loses its paragraph-level eastAsia rFonts hint,
# and four new entries (an empty line, "9.26", a css entry and a follow-up js
# note) are inserted after it, followed by a brand new closing paragraph (the
# one that now carries the old paragraph's <w:pPr> rFonts hint) comparing
# html/css understanding.
#
# Word's Range.Find/InsertAfter can't place <w:proofErr> spell-check markers
# around the "js"/"css" runs, so the new paragraphs are authored as literal
# WordprocessingML and dropped in via Range.InsertXML -- the COM-exposed
# equivalent of pasting OOXML into a Range.

$d = $word.ActiveDocument

# Paragraph 2 is the existing "9.24/9.25" journal entry (工欲善其事...).
$entry = $d.Paragraphs(2)

$newBody = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>工欲善其事，必先利其器。学习一些基本用法，笔记里写的就是我掌握的所有东西。探索GitHub的使用，还是有点懵逼，这几天的上传可能会有点混乱。试着写了一点代码。迫切的想要知道的是</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>js</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>的代码是怎样控制网页运行的。</w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>9</w:t></w:r><w:r><w:t>.26</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>在网页上调试，深刻的明白了什么是</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>css</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>万物皆盒，由此写出来一些东西。需要深入了解的是定位与布局，不然很难进步。</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>另有：知道了</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>js</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>怎么控制网页运行。准备尝试。</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>比较而言，我觉得</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>js</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>部分更好理解了，html和</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>css</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>看似简单，却是一个新的脑回路。</w:t></w:r></w:p>
'@

$entry.Range.InsertXML($newBody)
